$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.706.10'
$ws.Range('E2').Value = '  +3.78%  '
$ws.Range('D3').Value = '3.439.33'
$ws.Range('E3').Value = '  +2.85%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.60%  '
$ws.Range('E7').Value = '  +2.05%  '
$ws.Range('D8').Value = '3.434.07'
$ws.Range('E8').Value = '  +2.97%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.172'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.643'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.09'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.89%  '
$ws.Range('E13').Value = '  +0.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').Value = '3.986.33'
$ws.Range('E15').Value = '  +2.65%  '
$ws.Range('E16').Value = '  +2.40%  '
$ws.Range('D17').Value = '3.442.37'
$ws.Range('E17').Value = '  +2.74%  '
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '66.687.04'
$ws.Range('E19').Value = '  +2.40%  '
$ws.Range('E20').Value = '  +2.70%  '
$ws.Range('E21').Value = '  +2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '486.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +17.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.02'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.28'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.98%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.46%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.30%  '
$ws.Range('E29').Value = '  +5.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.06'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.98%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '593.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.52%  '
$ws.Range('E35').Value = '  +3.89%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.148'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('E39').Value = '  +4.47%  '
$ws.Range('E40').Value = '  +4.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.28'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.45%  '
$ws.Range('D42').Value = '3.183.18'
$ws.Range('E42').Value = '  +3.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.77%  '
$ws.Range('E44').Value = '  +3.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.53'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.87%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +21.54%  '
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.72'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '140.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.87%  '
